$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 84; this shifts old rows 84-88 down to 85-89
$ws.Rows.Item(84).Insert()

# Copy style from row 85 (the old row 84, now shifted down) column D onto new row 84 column D
$ws.Cells.Item(85, 4).Copy()
$ws.Cells.Item(84, 4).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Populate the new row 84 with data
$ws.Cells.Item(84, 1).Value = 1
$ws.Cells.Item(84, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(84, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(84, 4).Value = 45132
$ws.Cells.Item(84, 5).Value = 15
$ws.Cells.Item(84, 6).Value = 100112031
$ws.Cells.Item(84, 7).Value = "Poroto verde"
$ws.Cells.Item(84, 8).Value = "Magnum"
$ws.Cells.Item(84, 9).Value = "Primera"
$ws.Cells.Item(84, 10).Value = 160
$ws.Cells.Item(84, 11).Value = 19000
$ws.Cells.Item(84, 12).Value = 20000
$ws.Cells.Item(84, 13).Value = 19438
$ws.Cells.Item(84, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(84, 15).Value = "Perú"
$ws.Cells.Item(84, 16).Value = 778
$ws.Cells.Item(84, 17).Value = 25
$ws.Cells.Item(84, 18).Value = "Hortaliza"
